$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newDate = "2025-11-02 01:24:15"

# --- Insert a new row for the new listing, pushing the existing rows 3-5 down to 4-6 ---
$ws.Rows.Item(3).Insert()

# --- Refresh the fetch timestamp for every data row (2-6) ---
$ws.Range("A2").Value = $newDate
$ws.Range("A3").Value = $newDate
$ws.Range("A4").Value = $newDate
$ws.Range("A5").Value = $newDate
$ws.Range("A6").Value = $newDate

# --- Fill in the newly inserted row 3 with the new listing's data ---
$ws.Range("B3").Value = "【急募】Transformerベースのテキストエンコーダー経験者募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5425363"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("G3").Value = 25

# --- Rebuild all hyperlinks on column F so they point at the correct rows ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5425201") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5425363") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5425003") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5425263") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5424906") | Out-Null

# --- Re-apply the hyperlink style to the F column cells (Hyperlinks.Delete/Add resets formatting) ---
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"

# --- Widen column D slightly so the longer price string fits ---
$ws.Columns.Item(4).ColumnWidth = 31.17
